$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B19").Value = "< 10 in Manuf., Mining, Constr., Transport `n < 5 in other "
$ws.Range("B20").Value = "< 50 in Manufacturing, Mining, Constr., Transport `n< 10 in other "
$ws.Range("B21").Value = "< 300 emp. Manuf., Mining, Construction, Transport`n<200 Agric. `n<100 Others`n<50 Real State"
$ws.Range("C21").Value = "< 8 Billionlion won Manuf., `n< 3 Billionlion won Mining, Construction, Transport"
$ws.Range("D21").Value = "<30 Billionlion won Publication, Inf. and Communication, and others, `n<20 Billionlion won Agriculture and others, `n<10 Billionlion won Sewerage, waste management, and others, `n<5 Billionlion won Real State"
$ws.Range("B22").Value = ">= 300 emp. Manuf., Mining, Construction, Transport`n>=200 Agric.`n>=100 Others `n>=50 Real State"
$ws.Range("C22").Value = ">=8 Billionlion won Manuf., `n>=3 Billionlion won Mining, Construction, Transport"
$ws.Range("D22").Value = "Billionlion Won`n>30 Publication, Inf., Communication, and others, `n>20 Agriculture and others, `n>10 Sewerage, waste management, and others, `n>5 Real State"
